# Update the "Generate Report for Handback" timestamps that were refreshed
# when the report was regenerated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-02 19:13:15"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2) for the first file (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-02 19:13:09"
$wsZhCn.Range("K2").Value = "2016-09-02 19:13:37"

# de-de sheet: "Correspond Handback DateTime" (K2) for the first file (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-02 19:13:44"
